$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 23250
$ws.Range("D2").Value = -0.0292
$ws.Range("I2").Value = 4.3
$ws.Range("J2").Value = 79
$ws.Range("K2").Value = 79
$ws.Range("C3").Value = 106000
$ws.Range("D3").Value = 0.0047
$ws.Range("I3").Value = 6.13
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 70
$ws.Range("C4").Value = 439000
$ws.Range("D4").Value = 0.008
$ws.Range("D4").NumberFormat = "0.00%"
$ws.Range("I4").Value = 4.33
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 76
$ws.Range("C5").Value = 32200
$ws.Range("D5").Value = 0.0078
$ws.Range("I5").Value = 6.21
$ws.Range("J5").Value = 51
$ws.Range("K5").Value = 51
$ws.Range("C6").Value = 32700
$ws.Range("D6").Value = 0
$ws.Range("D6").NumberFormat = "0%"
$ws.Range("I6").Value = 3.67
$ws.Range("J6").Value = 86
$ws.Range("K6").Value = 86
$ws.Range("C7").Value = 25650
$ws.Range("D7").Value = 0.0039
$ws.Range("I7").Value = 4.68
$ws.Range("J7").Value = 71
$ws.Range("K7").Value = 71
$ws.Range("C8").Value = 10610
$ws.Range("D8").Value = 0.0114
$ws.Range("I8").Value = 4.85
$ws.Range("J8").Value = 85
$ws.Range("K8").Value = 85
$ws.Range("C9").Value = 83500
$ws.Range("D9").Value = 0.053
$ws.Range("I9").Value = 3.59
$ws.Range("J9").Value = 74
$ws.Range("K9").Value = 74
$ws.Range("C10").Value = 217000
$ws.Range("D10").Value = -0.0069
$ws.Range("I10").Value = 5.53
$ws.Range("J10").Value = 53
$ws.Range("K10").Value = 53
$ws.Range("C11").Value = 130200
$ws.Range("D11").Value = 0.0148
$ws.Range("I11").Value = 5.22
$ws.Range("J11").Value = 84
$ws.Range("K11").Value = 84
$ws.Range("C12").Value = 20100
$ws.Range("D12").Value = 0
$ws.Range("D12").NumberFormat = "0%"
$ws.Range("I12").Value = 4.73
$ws.Range("J12").Value = 78
$ws.Range("K12").Value = 78
$ws.Range("C13").Value = 70700
$ws.Range("D13").Value = 0.0071
$ws.Range("I13").Value = 4.95
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 80
$ws.Range("C14").Value = 54900
$ws.Range("D14").Value = -0.0054
$ws.Range("D14").NumberFormat = "0.00%"
$ws.Range("I14").Value = 6.45
$ws.Range("J14").Value = 69
$ws.Range("K14").Value = 69
$ws.Range("C15").Value = 81500
$ws.Range("D15").Value = 0.0062
$ws.Range("I15").Value = 6.75
$ws.Range("J15").Value = 85
$ws.Range("K15").Value = 85
$ws.Range("C16").Value = 19270
$ws.Range("D16").Value = 0.0142
$ws.Range("I16").Value = 5.53
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 79
$ws.Range("C17").Value = 49700
$ws.Range("D17").Value = 0.0164
$ws.Range("I17").Value = 5.63
$ws.Range("J17").Value = 69
$ws.Range("K17").Value = 69
$ws.Range("C18").Value = 20250
$ws.Range("D18").Value = 0.01
$ws.Range("I18").Value = 6.07
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = 38
$ws.Range("C19").Value = 54500
$ws.Range("D19").Value = -0.0145
$ws.Range("I19").Value = 3.67
$ws.Range("J19").Value = 87
$ws.Range("K19").Value = 87
$ws.Range("C20").Value = 14260
$ws.Range("D20").Value = -0.0021
$ws.Range("I20").Value = 4.56
$ws.Range("J20").Value = 71
$ws.Range("K20").Value = 71
$ws.Range("C21").Value = 130800
$ws.Range("D21").Value = -0.0068
$ws.Range("I21").Value = 4.13
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 80
$ws.Range("C22").Value = 43750
$ws.Range("D22").Value = -0.0191
$ws.Range("I22").Value = 3.33
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 50
$ws.Range("C23").Value = 68900
$ws.Range("D23").Value = 0.0253
$ws.Range("I23").Value = 3.13
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 90
$ws.Range("C24").Value = 49450
$ws.Range("D24").Value = -0.009
$ws.Range("I24").Value = 5.46
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 70
$ws.Range("C25").Value = 85600
$ws.Range("D25").Value = 0.0154
$ws.Range("I25").Value = 4.21
$ws.Range("J25").Value = 84
$ws.Range("K25").Value = 84
$ws.Range("C26").Value = 111200
$ws.Range("D26").Value = 0.0063
$ws.Range("I26").Value = 2.85
$ws.Range("J26").Value = 83
$ws.Range("K26").Value = 83
$ws.Range("C27").Value = 14470
$ws.Range("D27").Value = 0.0028
$ws.Range("I27").Value = 4.49
$ws.Range("J27").Value = 86
$ws.Range("K27").Value = 86
$ws.Range("C28").Value = 13850
$ws.Range("D28").Value = 0.0344
$ws.Range("I28").Value = 3.61
$ws.Range("J28").Value = 83
$ws.Range("K28").Value = 83
$ws.Range("C29").Value = 23150
$ws.Range("D29").Value = 0.0266
$ws.Range("I29").Value = 4.3
$ws.Range("J29").Value = 86
$ws.Range("K29").Value = 86
$ws.Range("C30").Value = 24700
$ws.Range("D30").Value = 0.002
$ws.Range("I30").Value = 4.86
$ws.Range("J30").Value = 87
$ws.Range("K30").Value = 87

$null = $ws.Range("A1").Select()
